# Auto-generated edit script: refreshes cached market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 830.46155
$ws.Range("I41").Value = 1187.125
$ws.Range("K41").Value = 1187.125
$ws.Range("M41").Value = -747.125

$ws.Range("H74").Value = 7060.143
$ws.Range("I74").Value = 6595.1665
$ws.Range("J74").Value = 9850
$ws.Range("K74").Value = 6595.1665
$ws.Range("L74").Value = 9850
$ws.Range("M74").Value = -5659.1665
$ws.Range("N74").Value = -11722

$ws.Range("H77").Value = 7060.143
$ws.Range("I77").Value = 6595.1665
$ws.Range("J77").Value = 9850
$ws.Range("K77").Value = 32975.8325
$ws.Range("L77").Value = 49250
$ws.Range("M77").Value = -28295.8325
$ws.Range("N77").Value = -58610

$ws.Range("H87").Value = 75997.8
$ws.Range("I87").Value = 40000
$ws.Range("J87").Value = 99996.336
$ws.Range("K87").Value = 40000
$ws.Range("L87").Value = 99996.336
$ws.Range("M87").Value = -38752
$ws.Range("N87").Value = -102492.336

$ws.Range("H88").Value = 1841.7391
$ws.Range("J88").Value = 1924.1578
$ws.Range("L88").Value = 1924.1578
$ws.Range("N88").Value = -2736.1578

$ws.Range("H90").Value = 75997.8
$ws.Range("I90").Value = 40000
$ws.Range("J90").Value = 99996.336
$ws.Range("K90").Value = 120000
$ws.Range("L90").Value = 299989.008
$ws.Range("M90").Value = -113760
$ws.Range("N90").Value = -312469.008

$ws.Range("H91").Value = 1841.7391
$ws.Range("J91").Value = 1924.1578
$ws.Range("L91").Value = 1924.1578
$ws.Range("N91").Value = -4732.1578

$ws.Range("H100").Value = 3968.2856
$ws.Range("I100").Value = 1320.3077
$ws.Range("K100").Value = 1320.3077
$ws.Range("M100").Value = -779.3077000000001

$ws.Range("H107").Value = 1087.2106
$ws.Range("I107").Value = 762
$ws.Range("K107").Value = 762
$ws.Range("M107").Value = 1158

$ws.Range("H113").Value = 5002.6924
$ws.Range("J113").Value = 2860.75
$ws.Range("L113").Value = 2860.75
$ws.Range("N113").Value = -9368.75

$ws.Range("H138").Value = 4989.357
$ws.Range("J138").Value = 6855.2
$ws.Range("L138").Value = 20565.6
$ws.Range("N138").Value = -30845.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 10008
$ws.Range("I17").Value = 10008
$ws.Range("K17").Value = 10008
$ws.Range("M17").Value = -9835

$ws.Range("H74").Value = 1677.091
$ws.Range("I74").Value = 1444.8
$ws.Range("K74").Value = 1444.8
$ws.Range("M74").Value = -570.8

$ws.Range("H77").Value = 1677.091
$ws.Range("I77").Value = 1444.8
$ws.Range("K77").Value = 7224
$ws.Range("M77").Value = -2856

$ws.Range("H132").Value = 5885265
$ws.Range("I132").Value = 3193.6924
$ws.Range("K132").Value = 9581.0772
$ws.Range("M132").Value = -7051.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1624.2727
$ws.Range("I22").Value = 1741.8
$ws.Range("J22").Value = 1526.3334
$ws.Range("K22").Value = 1741.8
$ws.Range("L22").Value = 1526.3334
$ws.Range("M22").Value = -1568.8
$ws.Range("N22").Value = -1872.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 426.55554
$ws.Range("I7").Value = 200.5
$ws.Range("J7").Value = 607.4
$ws.Range("K7").Value = 200.5
$ws.Range("L7").Value = 607.4
$ws.Range("M7").Value = -87.5
$ws.Range("N7").Value = -833.4

$ws.Range("H28").Value = 15082.333
$ws.Range("J28").Value = 15082.333
$ws.Range("L28").Value = 15082.333
$ws.Range("N28").Value = -15572.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1145.5714
$ws.Range("I5").Value = 1136.8889
$ws.Range("J5").Value = 1161.2
$ws.Range("K5").Value = 3410.6667
$ws.Range("L5").Value = 3483.6
$ws.Range("M5").Value = -3298.6667
$ws.Range("N5").Value = -3707.6

$ws.Range("H12").Value = 3962.4
$ws.Range("I12").Value = 23.666666
$ws.Range("J12").Value = 5650.4287
$ws.Range("K12").Value = 70.99999800000001
$ws.Range("L12").Value = 16951.2861
$ws.Range("M12").Value = 102.000002
$ws.Range("N12").Value = -17297.2861

$ws.Range("H82").Value = 13619.286
$ws.Range("I82").Value = 10400.6
$ws.Range("J82").Value = 21666
$ws.Range("K82").Value = 31201.8
$ws.Range("L82").Value = 64998
$ws.Range("M82").Value = -30795.8
$ws.Range("N82").Value = -65810

$ws.Range("H85").Value = 13619.286
$ws.Range("I85").Value = 10400.6
$ws.Range("J85").Value = 21666
$ws.Range("K85").Value = 31201.8
$ws.Range("L85").Value = 64998
$ws.Range("M85").Value = -29797.8
$ws.Range("N85").Value = -67806

$ws.Range("H117").Value = 6187
$ws.Range("I117").Value = 231.66667
$ws.Range("J117").Value = 12142.333
$ws.Range("K117").Value = 695.00001
$ws.Range("L117").Value = 36426.999
$ws.Range("M117").Value = 2746.99999
$ws.Range("N117").Value = -43310.999

$ws.Range("H118").Value = 6541.8
$ws.Range("I118").Value = 3565
$ws.Range("K118").Value = 10695
$ws.Range("M118").Value = -9452

$ws.Range("H129").Value = 1195576.6
$ws.Range("I129").Value = 3389.625
$ws.Range("K129").Value = 10168.875
$ws.Range("M129").Value = -5168.875

$ws.Range("H135").Value = 1145.5714
$ws.Range("I135").Value = 1136.8889
$ws.Range("J135").Value = 1161.2
$ws.Range("K135").Value = 10232.0001
$ws.Range("L135").Value = 10450.8
$ws.Range("M135").Value = -7697.000099999999
$ws.Range("N135").Value = -15520.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2162.4546
$ws.Range("I80").Value = 1561.8889
$ws.Range("J80").Value = 4865
$ws.Range("K80").Value = 1561.8889
$ws.Range("L80").Value = 4865
$ws.Range("M80").Value = -563.8888999999999
$ws.Range("N80").Value = -6861

$ws.Range("H83").Value = 2162.4546
$ws.Range("I83").Value = 1561.8889
$ws.Range("J83").Value = 4865
$ws.Range("K83").Value = 7809.4445
$ws.Range("L83").Value = 24325
$ws.Range("M83").Value = -2817.4445
$ws.Range("N83").Value = -34309

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 25000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 25000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 25000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -25226

$ws.Range("H22").Value = 1939
$ws.Range("I22").Value = 831.5
$ws.Range("K22").Value = 831.5
$ws.Range("M22").Value = -536.5

$ws.Range("H27").Value = 1939
$ws.Range("I27").Value = 831.5
$ws.Range("K27").Value = 831.5
$ws.Range("M27").Value = -724.5

$ws.Range("H28").Value = 25000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 25000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 25000
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -25464

$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 25000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -25214

$ws.Range("H42").Value = 52000
$ws.Range("I42").Value = 52000
$ws.Range("K42").Value = 52000
$ws.Range("M42").Value = -51437

$ws.Range("H49").Value = 52000
$ws.Range("I49").Value = 52000
$ws.Range("K49").Value = 52000
$ws.Range("M49").Value = -51853

$ws.Range("H61").Value = 2355.75
$ws.Range("I61").Value = 2355.75
$ws.Range("K61").Value = 2355.75
$ws.Range("M61").Value = -2153.75

$ws.Range("H113").Value = 2355.75
$ws.Range("I113").Value = 2355.75
$ws.Range("K113").Value = 2355.75
$ws.Range("M113").Value = -185.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10066.5
$ws.Range("I29").Value = 3800
$ws.Range("J29").Value = 16333
$ws.Range("K29").Value = 3800
$ws.Range("L29").Value = 16333
$ws.Range("M29").Value = -3510
$ws.Range("N29").Value = -16913

$ws.Range("H45").Value = 21416.666
$ws.Range("J45").Value = 8125
$ws.Range("L45").Value = 8125
$ws.Range("N45").Value = -9107
